# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp string (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 30 de Marzo de 2020 a las 22:20"

# --- Estados Unidos (row 4): refreshed totals ---
$ws.Range("B4").Value = 159689
$ws.Range("C4").Value = 16198
$ws.Range("D4").Value = 5220
$ws.Range("E4").Value = 151518
$ws.Range("F4").Value = 3402
$ws.Range("G4").Value = 368
$ws.Range("H4").Value = 2951

# --- Alemania (row 8): refreshed totals ---
$ws.Range("B8").Value = 66125
$ws.Range("C8").Value = 3690
$ws.Range("D8").Value = 13500
$ws.Range("E8").Value = 52009
$ws.Range("F8").Value = 1979
$ws.Range("G8").Value = 75
$ws.Range("H8").Value = 616

# --- India's case count overtakes Grecia & Islandia, so it moves up to
#     row 43, pushing Grecia to row 44 and Islandia to row 45. ---

# Row 43 now shows India with its refreshed totals
$ws.Range("A43").Value = "India"
$ws.Range("B43").Value = 1251
$ws.Range("C43").Value = 227
$ws.Range("D43").Value = 100
$ws.Range("E43").Value = 1119
$ws.Range("F43").Value = 0
$ws.Range("G43").Value = 5
$ws.Range("H43").Value = 32

# Row 44 now shows Grecia (previously at row 43) with its prior totals
$ws.Range("A44").Value = "Grecia"
$ws.Range("B44").Value = 1212
$ws.Range("C44").Value = 56
$ws.Range("D44").Value = 52
$ws.Range("E44").Value = 1117
$ws.Range("F44").Value = 72
$ws.Range("G44").Value = 4
$ws.Range("H44").Value = 43

# Row 45 now shows Islandia (previously at row 44) with its prior totals
$ws.Range("A45").Value = "Islandia"
$ws.Range("B45").Value = 1086
$ws.Range("C45").Value = 66
$ws.Range("D45").Value = 157
$ws.Range("E45").Value = 927
$ws.Range("F45").Value = 25
$ws.Range("G45").Value = 0
$ws.Range("H45").Value = 2
